$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "region"
$ws.Range("B1").Value = "percent"

# Data rows (23 countries)
$data = @(
    @("China", 20),
    @("Japan", 30),
    @("UK", 40),
    @("Germany", 99),
    @("Brazil", 80),
    @("Canada", 20),
    @("ANZ", 40),
    @("Benelux", 50),
    @("Italy", 30),
    @("France", 20),
    @("Iberia", 50),
    @("Mexico", 60),
    @("Thailand", 40),
    @("India", 20),
    @("Korea", 50),
    @("Saudi", 10),
    @("Gulf", 60),
    @("Chile", 44),
    @("Turkey", 22),
    @("Poland", 66),
    @("S Africa", 77),
    @("Indonesia", 44),
    @("Russia", 22)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Borders: A2 gets left+top (table top-left corner), A3:A24 get left only.
# Build A2's combined border first (left, then top), then apply the
# left-only border to the rest of the column -- this ordering lines up
# the resulting style indexes with A2 -> s=1, A3:A24 -> s=2.
$ws.Range("A2").Borders.Item(7).LineStyle = 1
$ws.Range("A2").Borders.Item(7).Weight = 2
$ws.Range("A2").Borders.Item(8).LineStyle = 1
$ws.Range("A2").Borders.Item(8).Weight = 2

$ws.Range("A3:A24").Borders.Item(7).LineStyle = 1
$ws.Range("A3:A24").Borders.Item(7).Weight = 2

# Update selection to match the target state
$ws.Range("B25").Select()
